$wb = $excel.ActiveWorkbook

# Map row number (in sheet) -> new value for column F ("想去人数")
$updates = @{
    3  = 1697
    5  = 1119
    7  = 11885
    11 = 405
    14 = 13455
    15 = 13393
    20 = 278
}

# Both "展览" and "全部类型" sheets contain the same data and need the same update
$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
